# Natmi following Dr Hou advice
# Update the LR-pair results table (Col9a3 -> Mag) with recalculated values
# and add the new rows produced by the revised analysis (FAPs/sCs combinations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Col9a3"
$row2[0,2] = "Mag"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 0.5236536666666667
$row2[0,7] = 1.570961
$row2[0,8] = 0.4357111753954365
$row2[0,9] = 0.4357111753954365
$row2[0,10] = 1
$row2[0,11] = 0.3333333333333333
$row2[0,12] = 0.09045733333333333
$row2[0,13] = 0.271372
$row2[0,14] = 0.0574491187297735
$row2[0,15] = 0.0574491187297735
$row2[0,16] = 0.04736831427688889
$row2[0,17] = 0.426314828492
$row2[0,18] = 0.02503122304718159
$row2[0,19] = 0.02503122304718159
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Col9a3"
$row3[0,2] = "Mag"
$row3[0,3] = "sCs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 0.5236536666666667
$row3[0,7] = 1.570961
$row3[0,8] = 0.4357111753954365
$row3[0,9] = 0.4357111753954365
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 1.484107
$row3[0,13] = 4.452321
$row3[0,14] = 0.9425508812702265
$row3[0,15] = 0.9425508812702265
$row3[0,16] = 0.7771580722756667
$row3[0,17] = 6.994422650481001
$row3[0,18] = 0.4106799523482548
$row3[0,19] = 0.4106799523482548
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "FAPs"
$row4[0,1] = "Col9a3"
$row4[0,2] = "Mag"
$row4[0,3] = "ECs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 0.42003
$row4[0,7] = 1.26009
$row4[0,8] = 0.34949008600725
$row4[0,9] = 0.34949008600725
$row4[0,10] = 1
$row4[0,11] = 0.3333333333333333
$row4[0,12] = 0.09045733333333333
$row4[0,13] = 0.271372
$row4[0,14] = 0.0574491187297735
$row4[0,15] = 0.0574491187297735
$row4[0,16] = 0.03799479372
$row4[0,17] = 0.34195314348
$row4[0,18] = 0.02007789744590925
$row4[0,19] = 0.02007789744590926
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Col9a3"
$row5[0,2] = "Mag"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 0.42003
$row5[0,7] = 1.26009
$row5[0,8] = 0.34949008600725
$row5[0,9] = 0.34949008600725
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 1.484107
$row5[0,13] = 4.452321
$row5[0,14] = 0.9425508812702265
$row5[0,15] = 0.9425508812702265
$row5[0,16] = 0.62336946321
$row5[0,17] = 5.61032516889
$row5[0,18] = 0.3294121885613407
$row5[0,19] = 0.3294121885613407
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "sCs"
$row6[0,1] = "Col9a3"
$row6[0,2] = "Mag"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 0.258153
$row6[0,7] = 0.774459
$row6[0,8] = 0.2147987385973136
$row6[0,9] = 0.2147987385973136
$row6[0,10] = 1
$row6[0,11] = 0.3333333333333333
$row6[0,12] = 0.09045733333333333
$row6[0,13] = 0.271372
$row6[0,14] = 0.0574491187297735
$row6[0,15] = 0.0574491187297735
$row6[0,16] = 0.023351831972
$row6[0,17] = 0.210166487748
$row6[0,18] = 0.01233999823668265
$row6[0,19] = 0.01233999823668265
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "sCs"
$row7[0,1] = "Col9a3"
$row7[0,2] = "Mag"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 0.258153
$row7[0,7] = 0.774459
$row7[0,8] = 0.2147987385973136
$row7[0,9] = 0.2147987385973136
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 1.484107
$row7[0,13] = 4.452321
$row7[0,14] = 0.9425508812702265
$row7[0,15] = 0.9425508812702265
$row7[0,16] = 0.3831266743710001
$row7[0,17] = 3.448140069339
$row7[0,18] = 0.2024587403606309
$row7[0,19] = 0.2024587403606309
$ws.Range("A7:T7").Value = $row7
